$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'wedding leggings',
    'wedding running',
    'wedding workout clothes',
    'weight lifting ankle support',
    'weight vest basketball',
    'weightlifting apparel',
    'weightlifting compression',
    'weightlifting floor',
    'weightlifting leggings women',
    'weightlifting pants women',
    'weightlifting shirts for women',
    'weightlifting shoes women black',
    'weightlifting tanks for women',
    'weightlifting tape',
    'weightlifting vest',
    'weird addictions',
    'weird things for women',
    'western riding gear',
    'which tights women',
    'white athletic tights',
    'white basketball leggings',
    'white capri leggins for women',
    'white capris for men',
    'white capris leggings',
    'white capris plus',
    'white capris women',
    'white cheerleader shorts',
    'white compression hose',
    'white compression leggings women',
    'white compression pants women',
    'white compression tights',
    'white compression tights boys',
    'white cycling shorts women',
    'white dress panta',
    'white football tights',
    'white full length leggings',
    'white knee high tights',
    'white knee patella',
    'white knee support',
    'white legging capris',
    'white leggings capri',
    'white leggings shorts',
    'white leggings with holes',
    'white lululemon leggings',
    'white lycra pants',
    'white lycra shorts',
    'white mens compression tights',
    'white running leggings for women',
    'white running pants women',
    'white running tights men',
    'white running tights women',
    'white skirt leggings',
    'white spandex capri leggings',
    'white spandex capris',
    'white tennis skirt lululemon',
    'white thigh highs stay up',
    'white tights capri',
    'white tights for men basketball',
    'white tights for teens',
    'white tights for women with design',
    'white tights mens basketball',
    'white weightlifting tape',
    'white yoga capri pants for women',
    'wick inhaler',
    'wicking running gear',
    'wide band tights for women',
    'wide hips',
    'wide leg yoga shorts',
    'wide waistband stripe side leggings',
    'wide waistband yoga',
    'wider hips',
    'wind running tights',
    'winter athletic gear',
    'winter athletic skirts for women',
    'winter athletic wear women',
    'winter ball dresses for teens cheap',
    'winter basketball',
    'winter bike pants women',
    'winter biking gear',
    'winter clothes for runners',
    'winter clthes for women',
    'winter compression',
    'winter cycling clothes women',
    'winter cycling tights women',
    'winter equipment',
    'winter legginga',
    'winter leggingd',
    'winter leggings',
    'winter leggings for plus women',
    'winter leggings for women',
    'winter leggings for womens',
    'winter leggings men',
    'winter leggins for women',
    'winter leggins for womens',
    'winter men leggings',
    'winter outdoor running gear women',
    'winter outdoor running pants women',
    'winter outdoor skirts for women',
    'winter overall',
    'winter panta'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
